$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D2, E2
$ws.Range('D2').Value = '66.341.90'
$ws.Range('E2').Value = '  +0.15%  '

# Row 3: update D3, E3
$ws.Range('D3').Value = '3.555.59'
$ws.Range('E3').Value = '  +0.72%  '

# Row 4: update E4
$ws.Range('E4').Value = '  -0.09%  '

# Row 5: update D5, E5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.41'
$ws.Range('E5').Value = '  -0.43%  '

# Row 6: update D6, E6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.32'
$ws.Range('E6').Value = '  +0.25%  '

# Row 7: update D7, E7
$ws.Range('D7').Value = '3.554.81'
$ws.Range('E7').Value = '  +0.78%  '

# Row 8: update E8
$ws.Range('E8').Value = '  +0.12%  '

# Row 9: update D9, E9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  +2.39%  '

# Row 10: update E10
$ws.Range('E10').Value = '  -0.66%  '

# Row 11: update E11
$ws.Range('E11').Value = '  -3.00%  '

# Row 12: update D12, E12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  +0.42%  '

# Row 13: update D13, E13
$ws.Range('D13').Value = '4.155.12'
$ws.Range('E13').Value = '  +0.60%  '

# Row 14: update D14, E14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000207'
$ws.Range('E14').Value = '  -0.11%  '

# Row 15: update D15, E15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '30.16'
$ws.Range('E15').Value = '  -0.47%  '

# Row 16: update D16, E16
$ws.Range('D16').Value = '3.551.71'
$ws.Range('E16').Value = '  +0.58%  '

# Row 17: update D17, E17
$ws.Range('D17').Value = '66.404.71'
$ws.Range('E17').Value = '  +0.14%  '

# Row 18: update E18
$ws.Range('E18').Value = '  +0.06%  '

# Row 19: update D19, E19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.48'
$ws.Range('E19').Value = '  +4.72%  '

# Row 20: update E20
$ws.Range('E20').Value = '  -0.57%  '

# Row 21: update D21, E21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.84'
$ws.Range('E21').Value = '  -0.69%  '

# Row 22: update D22, E22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '430.44'
$ws.Range('E22').Value = '  +1.08%  '

# Row 23: update D23, E23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.610'
$ws.Range('E23').Value = '  +1.37%  '

# Row 24: update D24, E24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.57'
$ws.Range('E24').Value = '  +1.02%  '

# Row 25: update D25, E25
$ws.Range('D25').Value = '3.698.53'
$ws.Range('E25').Value = '  +0.64%  '

# Row 26: update D26, E26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.09%  '

# Row 27: update E27
$ws.Range('E27').Value = '  -0.39%  '

# Row 28: update E28
$ws.Range('E28').Value = '  +1.35%  '

# Row 29: update D29, E29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  -1.09%  '

# Row 30: update D30, E30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.95'
$ws.Range('E30').Value = '  -0.85%  '

# Row 31: update D31, E31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.11%  '

# Row 32: update B32, C32, D32, E32
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.45'
$ws.Range('E32').Value = '  +0.56%  '

# Row 33: update B33, C33, D33, E33
$ws.Range('B33').Value = 'RenzoRestakedETH'
$ws.Range('C33').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D33').Value = '3.548.84'
$ws.Range('E33').Value = '  +0.75%  '

# Row 34: update D34, E34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.45'
$ws.Range('E34').Value = '  -2.05%  '

# Row 35: update E35
$ws.Range('E35').Value = '  -4.83%  '

# Row 36: update E36
$ws.Range('E36').Value = '  +0.01%  '

# Row 37: update D37, E37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.82'
$ws.Range('E37').Value = '  -0.31%  '

# Row 38: update E38
$ws.Range('E38').Value = '  -1.49%  '

# Row 39: update E39
$ws.Range('E39').Value = '  -0.42%  '

# Row 40: update D40, E40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '176.03'
$ws.Range('E40').Value = '  +1.78%  '

# Row 41: update D41, E41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0849'
$ws.Range('E41').Value = '  -0.76%  '

# Row 42: update D42, E42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.19'
$ws.Range('E42').Value = '  +0.26%  '

# Row 43: update D43, E43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.889'
$ws.Range('E43').Value = '  -0.35%  '

# Row 44: update D44, E44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.92'
$ws.Range('E44').Value = '  +1.40%  '

# Row 45: update D45, E45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '46.04'
$ws.Range('E45').Value = '  +1.66%  '

# Row 46: update D46, E46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.09%  '

# Row 47: update E47
$ws.Range('E47').Value = '  +4.00%  '

# Row 48: update E48
$ws.Range('E48').Value = '  -1.60%  '

# Row 49: update D49, E49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.10'
$ws.Range('E49').Value = '  -3.46%  '

# Row 50: update D50, E50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.14'
$ws.Range('E50').Value = '  -0.58%  '

# Row 51: update D51, E51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.19'
$ws.Range('E51').Value = '  +2.11%  '
